$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.962.77"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "2.933.96"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'374.32"
$ws.Range("E5").Value = "  -1.62%  "

$ws.Range("D6").Value = "'101.68"
$ws.Range("E6").Value = "  -2.76%  "

$ws.Range("E7").Value = "  -0.90%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "'0.581"
$ws.Range("E9").Value = "  -1.54%  "

$ws.Range("D10").Value = "'36.29"
$ws.Range("E10").Value = "  -1.71%  "

$ws.Range("E11").Value = "  -0.59%  "

$ws.Range("D12").Value = "'0.0835"
$ws.Range("E12").Value = "  -0.47%  "

$ws.Range("D13").Value = "3.397.00"
$ws.Range("E13").Value = "  -0.53%  "

$ws.Range("D14").Value = "'17.92"
$ws.Range("E14").Value = "  -2.40%  "

$ws.Range("D15").Value = "'7.32"
$ws.Range("E15").Value = "  -1.65%  "

$ws.Range("D16").Value = "2.912.32"
$ws.Range("E16").Value = "  -1.15%  "

$ws.Range("D17").Value = "'0.975"
$ws.Range("E17").Value = "  +1.96%  "

$ws.Range("D18").Value = "50.948.16"
$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("E19").Value = "  -5.64%  "

$ws.Range("D20").Value = "'7.15"
$ws.Range("E20").Value = "  -2.45%  "

$ws.Range("D21").Value = "'12.50"
$ws.Range("E21").Value = "  -2.70%  "

$ws.Range("D22").Value = "0.0₃0954"
$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").Value = "'264.49"
$ws.Range("E23").Value = "  +1.66%  "

$ws.Range("D24").Value = "'68.27"
$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("D25").Value = "'2.88"
$ws.Range("E25").Value = "  +2.52%  "

$ws.Range("D26").Value = "'8.06"
$ws.Range("E26").Value = "  +8.25%  "

$ws.Range("D27").Value = "'7.74"
$ws.Range("E27").Value = "  +8.52%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  -0.73%  "

$ws.Range("E30").Value = "  -0.23%  "

$ws.Range("D31").Value = "'25.58"
$ws.Range("E31").Value = "  -0.85%  "

$ws.Range("D32").Value = "'9.84"
$ws.Range("E32").Value = "  +0.43%  "

$ws.Range("D33").Value = "'50.96"
$ws.Range("E33").Value = "  -0.26%  "

$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").Value = "'0.0450"
$ws.Range("E34").Value = "  +0.83%  "

$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'33.51"
$ws.Range("E35").Value = "  -2.96%  "

$ws.Range("E36").Value = "  -3.04%  "

$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("D38").Value = "'2.97"
$ws.Range("E38").Value = "  -2.85%  "

$ws.Range("D39").Value = "'2.54"
$ws.Range("E39").Value = "  -1.32%  "

$ws.Range("E40").Value = "  -0.87%  "

$ws.Range("D41").Value = "'16.35"
$ws.Range("E41").Value = "  -5.06%  "

$ws.Range("E42").Value = "  -2.82%  "

$ws.Range("D43").Value = "'120.52"
$ws.Range("E43").Value = "  -1.97%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.285"
$ws.Range("E44").Value = "  +2.17%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'20.94"
$ws.Range("E45").Value = "  -4.52%  "

$ws.Range("D46").Value = "'2.05"
$ws.Range("E46").Value = "  -1.36%  "

$ws.Range("E47").Value = "  -1.71%  "

$ws.Range("D48").Value = "'3.21"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("D49").Value = "1.970.63"
$ws.Range("E49").Value = "  -3.05%  "

$ws.Range("D50").Value = "'0.0342"
$ws.Range("E50").Value = "  -1.61%  "

$ws.Range("D51").Value = "'5.02"
$ws.Range("E51").Value = "  -1.56%  "
